# Remove the "pander(table_forecasts)" SourceCode paragraph that
# immediately followed the "Using historical data..." paragraph and
# preceded the forecasts table.
$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*pander(table_forecasts)*") {
        $p.Range.Delete()
        break
    }
}
